$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-17: full Sending/Target cluster grid (ECs, FAPs, M2, sCs) x
# (ECs, FAPs, M2, sCs) for the Sema6a -> Plxna2 pair, recomputed per Dr Hou advice.
$data = @(
  @("ECs", "Sema6a", "Plxna2", "ECs", 3, 1, 53.11151633333333, 159.334549, 0.7183308303902657, 0.7183308303902658, 3, 1, 42.23928166666666, 126.717845, 0.7297675404946526, 0.7297675404946528, 2243.392298147434, 20190.53068332691, 0.5242145233553857, 0.5242145233553859),
  @("ECs", "Sema6a", "Plxna2", "FAPs", 3, 1, 53.11151633333333, 159.334549, 0.7183308303902657, 0.7183308303902658, 3, 1, 4.624984, 13.874952, 0.07990579066051323, 0.07990579066051323, 245.6399132574054, 2210.759219316648, 0.0573987929581572, 0.05739879295815722),
  @("ECs", "Sema6a", "Plxna2", "M2", 3, 1, 53.11151633333333, 159.334549, 0.7183308303902657, 0.7183308303902658, 3, 1, 2.151907, 6.455721, 0.03717847029587412, 0.03717847029587412, 114.2910437783143, 1028.619394004829, 0.02670644144027508, 0.02670644144027508),
  @("ECs", "Sema6a", "Plxna2", "sCs", 3, 1, 53.11151633333333, 159.334549, 0.7183308303902657, 0.7183308303902658, 3, 1, 8.864288333333333, 26.592865, 0.1531481985489599, 0.15314819854896, 470.7957945992094, 4237.162151392885, 0.1100110726364477, 0.1100110726364477),
  @("FAPs", "Sema6a", "Plxna2", "ECs", 3, 1, 4.187132, 12.561396, 0.05663076888328195, 0.05663076888328195, 3, 1, 42.23928166666666, 126.717845, 0.7297675404946526, 0.7297675404946528, 176.8614479235133, 1591.75303131162, 0.04132729692427378, 0.04132729692427378),
  @("FAPs", "Sema6a", "Plxna2", "FAPs", 3, 1, 4.187132, 12.561396, 0.05663076888328195, 0.05663076888328195, 3, 1, 4.624984, 13.874952, 0.07990579066051323, 0.07990579066051323, 19.365418505888, 174.288766552992, 0.004525126363331433, 0.004525126363331433),
  @("FAPs", "Sema6a", "Plxna2", "M2", 3, 1, 4.187132, 12.561396, 0.05663076888328195, 0.05663076888328195, 3, 1, 2.151907, 6.455721, 0.03717847029587412, 0.03717847029587412, 9.010318660724, 81.09286794651601, 0.00210544535875961, 0.00210544535875961),
  @("FAPs", "Sema6a", "Plxna2", "sCs", 3, 1, 4.187132, 12.561396, 0.05663076888328195, 0.05663076888328195, 3, 1, 8.864288333333333, 26.592865, 0.1531481985489599, 0.15314819854896, 37.11594533772666, 334.04350803954, 0.008672900236917126, 0.008672900236917127),
  @("M2", "Sema6a", "Plxna2", "ECs", 1, 0.3333333333333333, 0.02264733333333334, 0.067942, 0.0003063041479997877, 0.0003063041479997877, 3, 1, 42.23928166666666, 126.717845, 0.7297675404946526, 0.7297675404946528, 0.9566070916655556, 8.60946382499, 0.0002235308247291152, 0.0002235308247291152),
  @("M2", "Sema6a", "Plxna2", "FAPs", 1, 0.3333333333333333, 0.02264733333333334, 0.067942, 0.0003063041479997877, 0.0003063041479997877, 3, 1, 4.624984, 13.874952, 0.07990579066051323, 0.07990579066051323, 0.1047435543093334, 0.9426919887840001, 0.0000244754751285179, 0.0000244754751285179),
  @("M2", "Sema6a", "Plxna2", "M2", 1, 0.3333333333333333, 0.02264733333333334, 0.067942, 0.0003063041479997877, 0.0003063041479997877, 3, 1, 2.151907, 6.455721, 0.03717847029587412, 0.03717847029587412, 0.04873495513133334, 0.438614596182, 0.00001138791966791314, 0.00001138791966791314),
  @("M2", "Sema6a", "Plxna2", "sCs", 1, 0.3333333333333333, 0.02264733333333334, 0.067942, 0.0003063041479997877, 0.0003063041479997877, 3, 1, 8.864288333333333, 26.592865, 0.1531481985489599, 0.15314819854896, 0.2007524926477778, 1.80677243383, 0.00004690992847424151, 0.00004690992847424152),
  @("sCs", "Sema6a", "Plxna2", "ECs", 3, 1, 16.61610766666667, 49.848323, 0.2247320965784526, 0.2247320965784526, 3, 1, 42.23928166666666, 126.717845, 0.7297675404946526, 0.7297675404946528, 701.8524519359928, 6316.672067423935, 0.1640021893902641, 0.1640021893902641),
  @("sCs", "Sema6a", "Plxna2", "FAPs", 3, 1, 16.61610766666667, 49.848323, 0.2247320965784526, 0.2247320965784526, 3, 1, 4.624984, 13.874952, 0.07990579066051323, 0.07990579066051323, 76.84923210061068, 691.643088905496, 0.01795739586389607, 0.01795739586389607),
  @("sCs", "Sema6a", "Plxna2", "M2", 3, 1, 16.61610766666667, 49.848323, 0.2247320965784526, 0.2247320965784526, 3, 1, 2.151907, 6.455721, 0.03717847029587412, 0.03717847029587412, 35.75631840065367, 321.806865605883, 0.008355195577171511, 0.008355195577171511),
  @("sCs", "Sema6a", "Plxna2", "sCs", 3, 1, 16.61610766666667, 49.848323, 0.2247320965784526, 0.2247320965784526, 3, 1, 8.864288333333333, 26.592865, 0.1531481985489599, 0.15314819854896, 147.2899693350439, 1325.609724015395, 0.0344173157471209, 0.0344173157471209)
)

$r = 2
foreach ($row in $data) {
    for ($c = 1; $c -le $row.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    $r++
}
